# Adds the "CRIMINAL RECORD" block to the bookmarked paragraph (right after
# the {#show_photo}...{/show_photo} paragraph) and gives the
# {#show_border}...{/show_border} block its own paragraph, per the commit
# "add photo and border checkbox logic, fix-5".
$d = $word.ActiveDocument

# Locate the two template paragraphs we need to restructure.
$photoPara = $null
$borderPara = $null
foreach ($para in $d.Paragraphs) {
    $text = $para.Range.Text
    if ($text -like "*show_photo*") {
        $photoPara = $para
    }
    if ($text -like "*show_border*") {
        $borderPara = $para
    }
}

if ($photoPara -eq $null -or $borderPara -eq $null) {
    throw "Could not locate the show_photo / show_border template paragraphs"
}

# Replace the whole span (show_photo paragraph through show_border paragraph,
# including its trailing paragraph mark) with the restructured OOXML:
#   1) {#show_photo}PHOTO: in attachment{/show_photo}            (lang -> en-GB)
#   2) _GoBack bookmark + "CRIMINAL RECORD: {criminal_records}"  (new paragraph)
#   3) {#show_border}BORDER RECORD: ...{/show_border}            (new paragraph)
$target = $d.Range($photoPara.Range.Start, $borderPara.Range.End)

$newXml = @'
<w:p w:rsidR="00340EB6" w:rsidRPr="00790FC3" w:rsidRDefault="00C663B4" w:rsidP="00340EB6"><w:pPr><w:jc w:val="both"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00C663B4"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>{#</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00C663B4"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>show_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00C663B4"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>photo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00C663B4"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>}PHOTO</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00C663B4"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>: in attachment{/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00C663B4"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>show_photo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00C663B4"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>}</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="00340EB6"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">CRIMINAL RECORD: </w:t></w:r><w:r w:rsidR="005C2EC9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="005C2EC9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>criminal_</w:t></w:r><w:r w:rsidR="00790FC3"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>records</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00790FC3"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00C663B4"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>{#</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00C663B4"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>show_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00C663B4"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>border</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00C663B4"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>}BORDER</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00C663B4"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> RECORD: in attachment with the legend for translation{/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00C663B4"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>show_border</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00C663B4"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>}</w:t></w:r></w:p>
'@

$target.InsertXML($newXml) | Out-Null
